$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new "Nr" column
$ws.Range("F1").Value = "Nr"

# Numeric id values for rows 2..12
$ws.Range("F2").Value = 10001
$ws.Range("F3").Value = 10002
$ws.Range("F4").Value = 10003
$ws.Range("F5").Value = 10004
$ws.Range("F6").Value = 10005
$ws.Range("F7").Value = 10006
$ws.Range("F8").Value = 10007
$ws.Range("F9").Value = 10008
$ws.Range("F10").Value = 10009
$ws.Range("F11").Value = 10010
$ws.Range("F12").Value = 10011

# Update the view / selection to match the saved state
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("L10").Select()
